# Update Sheet1 values for rows 2-25 across columns B, C, D, E, F, H, I, L, O
# (case with 380 kV results - pl_mw.xlsx row data updated)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cols = @("B", "C", "D", "E", "F", "H", "I", "L", "O")

$data = @(
    @(1.353609635285011, 0.2539931740461157, 0.0262438978229298, 0.08598567600448526, 0.6471769034445316, 0.07973214163530429, 0.5584063081531312, 0.2280261146700582, 2.192062247822491),
    @(1.21133559458292, 0.2344638234132219, 0.02439403277229957, 0.08708918564973622, 0.6477343778755795, 0.07973214163530429, 0.5688561124922735, 0.2173551475487017, 2.208821801407936),
    @(1.123850794926113, 0.2224143719674885, 0.02325291009247366, 0.08782035368923857, 0.6486532643167493, 0.07973214163530429, 0.5757531420222719, 0.2109065978476821, 2.221066464059561),
    @(1.088170069107662, 0.2174897871213091, 0.02278659348611711, 0.08813178773660102, 0.6491724050279259, 0.07973214163530429, 0.5786843933572214, 0.208304865413723, 2.226546576890811),
    @(1.082243568007698, 0.2166712077245165, 0.02270908442068986, 0.08818431501456825, 0.6492673378069966, 0.07973214163530429, 0.5791784054669655, 0.2078744289044749, 2.227486127075281),
    @(1.123369710529744, 0.2223480148875012, 0.02324662640356223, 0.08782449923339097, 0.6486596802007867, 0.07973214163530429, 0.5757921857016051, 0.2108714041644504, 2.2211383869749),
    @(1.30458138523602, 0.2472717602742591, 0.02560718275345408, 0.08635503779960008, 0.6472492646653905, 0.07973214163530429, 0.5619094616711777, 0.2243253383031316, 2.197434682087248),
    @(1.658842196621038, 0.2956717304811605, 0.0301929920555466, 0.08389892303556756, 0.6490731178696336, 0.07973214163530429, 0.5385112720372902, 0.2515272424180779, 2.16650863173686),
    @(1.918370087783728, 0.3309281768760854, 0.0335345537003704, 0.0823539801019777, 0.6532328283311841, 0.07973214163530429, 0.5236673728869761, 0.2720110635473674, 2.153343904008125),
    @(2.036258689377291, 0.3468989227060035, 0.03504847107449649, 0.0817075449238569, 0.6557422041477565, 0.07973214163530429, 0.5174271531601313, 0.2814379857905749, 2.149445806838713),
    @(2.080873502717623, 0.3529366265009628, 0.03562083679938866, 0.08147086642692081, 0.6567815295646113, 0.07973214163530429, 0.5151380844476279, 0.2850233025403099, 2.148271556569142),
    @(2.071266143728508, 0.3516367524358088, 0.03549760913795552, 0.08152147846115021, 0.6565537250154136, 0.07973214163530429, 0.5156277819944535, 0.2842504504359766, 2.148511006604309),
    @(2.039929732064138, 0.3473958520712586, 0.03509557861629276, 0.08168791061573089, 0.6558259228511645, 0.07973214163530429, 0.5172373460159214, 0.2817326409034138, 2.149343143416644),
    @(2.020731685120552, 0.3447968562188635, 0.03484920225260879, 0.08179091180541143, 0.6553917338313084, 0.07973214163530429, 0.5182328918215049, 0.280192431760014, 2.149892201169678),
    @(1.910662131981496, 0.3298830608551953, 0.0334354884657273, 0.08239736100694017, 0.6530812815730513, 0.07973214163530429, 0.5240855206254267, 0.2713971722760959, 2.153640828946067),
    @(1.843092332736603, 0.3207163678340521, 0.03256661488063628, 0.08278383987984128, 0.6518221918253104, 0.07973214163530429, 0.5278073509043821, 0.2660293499350104, 2.156476834326838),
    @(1.804211924772517, 0.315437587505528, 0.03206628248331356, 0.08301143743215533, 0.6511560611011404, 0.07973214163530429, 0.52999625607303, 0.2629521599109097, 2.158304732819545),
    @(1.791045001665054, 0.3136492049095807, 0.03189678007411345, 0.08308940900076855, 0.6509404828864334, 0.07973214163530429, 0.5307456524000003, 0.2619120375999984, 2.158957373690896),
    @(1.850286930839275, 0.3216928367809544, 0.03265916815858816, 0.08274214943284264, 0.6519502121554979, 0.07973214163530429, 0.5274061643683225, 0.2665997046955084, 2.156154569603331),
    @(2.049134748026518, 0.3486417834200211, 0.03521369001480679, 0.08163880523209066, 0.6560372756996316, 0.07973214163530429, 0.5167625682729238, 0.2824717613933245, 2.149090521954264),
    @(2.178934613860577, 0.3661956155600876, 0.03687782855263322, 0.08096499492720532, 0.6592277865956362, 0.07973214163530429, 0.510237634935816, 0.2929356403277836, 2.146233889677177),
    @(2.109673261190096, 0.3568323061826106, 0.03599015070216893, 0.08132029094269555, 0.6574773157259557, 0.07973214163530429, 0.5136805567428553, 0.2873426131604191, 2.147597064933763),
    @(1.847034357619293, 0.321251402327789, 0.03261732732742217, 0.08276098084614603, 0.6518921543430238, 0.07973214163530429, 0.5275873876430595, 0.2663418197601857, 2.156299650556178),
    @(1.563130837478695, 0.2826305307332575, 0.02895717090301275, 0.08451779714560992, 0.6480861064781536, 0.07973214163530429, 0.5444300298351052, 0.2440808119954312, 2.173201939448944),
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $rowValues[$j]
    }
}
